$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move H4/I4 data down to a new row 5 ---
$h4 = $ws.Range("H4").Value2
$i4 = $ws.Range("I4").Value2
$ws.Range("H4").Clear()
$ws.Range("I4").Clear()
$ws.Range("H5").Value = $h4
$ws.Range("I5").Value = $i4
$ws.Rows.Item(5).RowHeight = 19.5

# --- Update C9 / D9 formulas to anchor the row with $ ---
$ws.Range("C9").Formula = "=H$4/1000*A3"
$ws.Range("D9").Formula = "=I$4/1000*B3"

# --- New column widths (G:I) ---
$ws.Columns.Item(7).ColumnWidth = 19.28515625
$ws.Columns.Item(8).ColumnWidth = 12.85546875
$ws.Columns.Item(9).ColumnWidth = 18.42578125

# --- New header row 11 (order of assignment matters for shared-string indices) ---
$ws.Range("H11").Value = "money limit"
$ws.Range("G13").Value = "total price per token"
$ws.Range("G11").Value = "tk output rate"
$ws.Range("I11").Value = "max amount of input tokens"
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("I11").VerticalAlignment = -4108
$ws.Range("I11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 30

# --- New calc cells around rows 12/14 ---
$ws.Range("G12").Formula = "=(G4)/F4"
$ws.Range("H12").Value = 2
$ws.Range("I12").Formula = "=H12/G14*1000"
$ws.Range("G14").Formula = "=A3+B3*G12"
$ws.Rows.Item(14).RowHeight = 19.5

# --- Window/view bookkeeping to mirror the saved cursor + window geometry ---
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = 10215
$win.Top = 960
$win.Width = 18330
$win.Height = 15090
$ws.Range("I13").Select()
